# Refresh the crypto price/volume columns (D = Price, E = Volume(1h)) for
# rows 2-51, matching the latest GitHub Actions scrape.
#
# The sheet stores these as plain text cells (e.g. "26.370.93", "0.9966",
# "  +6.50%  "). Several of the new Price values are plain decimals
# ("0.9966", "334.46", ...) that Excel's COM layer would otherwise coerce
# to Number on assignment, so we detect that case and use the classic
# leading-apostrophe trick to force text - then reset .Style back to
# "Normal" so we don't leave a stray quote-prefix format behind on cells
# that were unstyled before the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, [string]$val) {
    if ($val -match '^[+-]?\d+(\.\d+)?$') {
        $cell.Value = "'" + $val
        $cell.Style = "Normal"
    } else {
        $cell.Value = $val
    }
}

# row, new Price (D), new Volume(1h) (E)
$data = @(
    @(2,  "26.370.93",   "  +6.50%  "),
    @(3,  "1.725.10",    "  +4.16%  "),
    @(4,  "0.9966",      "  -0.33%  "),
    @(5,  "334.46",      "  +5.45%  "),
    @(6,  "0.9969",      "  -0.16%  "),
    @(7,  "0.3710",      "  +2.17%  "),
    @(8,  "49.41",       "  +5.60%  "),
    @(9,  "0.3358",      "  +2.95%  "),
    @(10, "1.196",       "  +5.02%  "),
    @(11, "0.07480",     "  +5.86%  "),
    @(12, "0.9979",      "  +0.07%  "),
    @(13, "6.349",       "  +5.12%  "),
    @(14, "20.24",       "  +3.37%  "),
    @(15, "6.962",       "  +5.11%  "),
    @(16, "1.716.90",    "  +3.49%  "),
    @(17, "0.00001082",  "  +3.19%  "),
    @(18, "0.06667",     "  +0.68%  "),
    @(19, "82.30",       "  +4.19%  "),
    @(20, "0.9969",      "  -0.08%  "),
    @(21, "16.50",       "  +4.73%  "),
    @(22, "6.116",       "  +3.42%  "),
    @(23, "13.09",       "  +4.00%  "),
    @(24, "26.251.53",   "  +6.25%  "),
    @(25, "2.463",       "  +1.41%  "),
    @(26, "2.488",       "  +3.99%  "),
    @(27, "151.68",      "  +1.81%  "),
    @(28, "1.380",       "  +13.14%  "),
    @(29, "19.37",       "  +4.02%  "),
    @(30, "1.912.51",    "  +3.76%  "),
    @(31, "129.97",      "  +3.44%  "),
    @(32, "4.117",       "  +0.90%  "),
    @(33, "6.026",       "  +3.08%  "),
    @(34, "0.08583",     "  +1.70%  "),
    @(35, "1.706",       "  +2.45%  "),
    @(36, "13.06",       "  +6.12%  "),
    @(37, "5.420",       "  +3.83%  "),
    @(38, "0.02347",     "  +4.87%  "),
    @(39, "0.06283",     "  +4.05%  "),
    @(40, "8.647",       "  +5.42%  "),
    @(41, "0.2151",      "  +3.91%  "),
    @(42, "1.237",       "  -3.23%  "),
    @(43, "0.6226",      "  +5.07%  "),
    @(44, "14.46",       "  +13.52%  "),
    @(45, "0.9973",      "  -0.11%  "),
    @(46, "3.871",       "  +1.37%  "),
    @(47, "0.5960",      "  +5.65%  "),
    @(48, "129.15",      "  +3.13%  "),
    @(49, "2.035",       "  +4.34%  "),
    @(50, "0.07298",     "  +4.66%  "),
    @(51, "77.25",       "  +3.47%  ")
)

foreach ($row in $data) {
    $r = $row[0]
    Set-TextValue $ws.Cells.Item($r, 4) $row[1]
    Set-TextValue $ws.Cells.Item($r, 5) $row[2]
}
